# Apply updated cryptocurrency price/volume figures to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.481.06"
$ws.Range("E2").Value = "  +0.98%  "

$ws.Range("D3").Value = "1.725.80"
$ws.Range("E3").Value = "  +0.49%  "

$ws.Range("D4").Value = "'0.9994"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'244.89"
$ws.Range("E5").Value = "  +2.42%  "

$ws.Range("D6").Value = "'0.9998"
$ws.Range("E6").Value = "  -0.09%  "

$ws.Range("D7").Value = "'0.4804"
$ws.Range("E7").Value = "  +1.62%  "

$ws.Range("D8").Value = "'0.2682"
$ws.Range("E8").Value = "  +2.34%  "

$ws.Range("D9").Value = "'0.06183"
$ws.Range("E9").Value = "  -0.27%  "

$ws.Range("D10").Value = "1.730.47"
$ws.Range("E10").Value = "  +0.80%  "

$ws.Range("D11").Value = "'0.07169"
$ws.Range("E11").Value = "  +1.17%  "

$ws.Range("D12").Value = "'15.57"
$ws.Range("E12").Value = "  +1.98%  "

$ws.Range("D13").Value = "'0.6088"
$ws.Range("E13").Value = "  +3.04%  "

$ws.Range("D14").Value = "'4.519"
$ws.Range("E14").Value = "  +2.45%  "

$ws.Range("D15").Value = "'76.99"
$ws.Range("E15").Value = "  +1.29%  "

$ws.Range("D16").Value = "'0.9996"
$ws.Range("E16").Value = "  -0.10%  "

$ws.Range("D17").Value = "26.505.04"
$ws.Range("E17").Value = "  +1.07%  "

$ws.Range("D18").Value = "'0.9997"
$ws.Range("E18").Value = "  -0.09%  "

$ws.Range("D19").Value = "'0.000006948"
$ws.Range("E19").Value = "  +2.00%  "

$ws.Range("D20").Value = "'11.50"
$ws.Range("E20").Value = "  -0.17%  "

$ws.Range("D21").Value = "1.953.91"
$ws.Range("E21").Value = "  +0.97%  "

$ws.Range("D22").Value = "'4.509"
$ws.Range("E22").Value = "  -1.15%  "

$ws.Range("D23").Value = "'8.791"
$ws.Range("E23").Value = "  +0.85%  "

$ws.Range("D24").Value = "'5.252"
$ws.Range("E24").Value = "  -0.39%  "

$ws.Range("D25").Value = "'137.01"
$ws.Range("E25").Value = "  +1.50%  "

$ws.Range("D26").Value = "'15.31"
$ws.Range("E26").Value = "  +0.96%  "

$ws.Range("D27").Value = "'1.778"
$ws.Range("E27").Value = "  +1.31%  "

$ws.Range("E28").Value = "  -0.43%  "

$ws.Range("D29").Value = "'106.53"
$ws.Range("E29").Value = "  -0.60%  "

$ws.Range("D30").Value = "'3.970"
$ws.Range("E30").Value = "  +0.20%  "

$ws.Range("D31").Value = "'0.08002"
$ws.Range("E31").Value = "  +3.31%  "

$ws.Range("D32").Value = "'3.689"
$ws.Range("E32").Value = "  +0.43%  "

$ws.Range("D33").Value = "'0.04520"
$ws.Range("E33").Value = "  +1.97%  "

$ws.Range("D34").Value = "'0.9993"
$ws.Range("E34").Value = "  -0.09%  "

$ws.Range("D35").Value = "'2.614"
$ws.Range("E35").Value = "  -0.06%  "

$ws.Range("D36").Value = "'0.9904"
$ws.Range("E36").Value = "  +1.89%  "

$ws.Range("D37").Value = "'0.6284"
$ws.Range("E37").Value = "  +2.13%  "

$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'0.9130"
$ws.Range("E38").Value = "  -1.23%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "'2.079"
$ws.Range("E39").Value = "  +8.47%  "

$ws.Range("D40").Value = "'2.375"
$ws.Range("E40").Value = "  -2.61%  "

$ws.Range("D41").Value = "'1.002"
$ws.Range("E41").Value = "  +0.20%  "

$ws.Range("D42").Value = "'103.92"
$ws.Range("E42").Value = "  -7.47%  "

$ws.Range("D43").Value = "'0.01501"
$ws.Range("E43").Value = "  +1.98%  "

$ws.Range("D44").Value = "'5.607"
$ws.Range("E44").Value = "  +4.74%  "

$ws.Range("D45").Value = "'0.3866"
$ws.Range("E45").Value = "  +1.56%  "

$ws.Range("D46").Value = "'6.913"
$ws.Range("E46").Value = "  +10.42%  "

$ws.Range("D47").Value = "'0.1180"
$ws.Range("E47").Value = "  +1.07%  "

$ws.Range("E48").Value = "  +1.57%  "

$ws.Range("D49").Value = "'30.53"
$ws.Range("E49").Value = "  +1.31%  "

$ws.Range("D50").Value = "'7.780"
$ws.Range("E50").Value = "  +0.98%  "

$ws.Range("E51").Value = "  +3.23%  "
